$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.71508866666667
$ws.Range("H2").Value = 107.145266
$ws.Range("I2").Value = 0.1390302752364672
$ws.Range("J2").Value = 0.1390302752364672
$ws.Range("M2").Value = 0.9785706666666668
$ws.Range("N2").Value = 2.935712000000001
$ws.Range("O2").Value = 0.04994923768046061
$ws.Range("P2").Value = 0.04994923768046061
$ws.Range("Q2").Value = 34.94973812659912
$ws.Range("R2").Value = 314.5476431393921
$ws.Range("S2").Value = 0.006944456262566155
$ws.Range("T2").Value = 0.006944456262566156
$ws.Range("G3").Value = 35.71508866666667
$ws.Range("H3").Value = 107.145266
$ws.Range("I3").Value = 0.1390302752364672
$ws.Range("J3").Value = 0.1390302752364672
$ws.Range("O3").Value = 0.5521442932757068
$ws.Range("P3").Value = 0.5521442932757068
$ws.Range("Q3").Value = 386.3381976224016
$ws.Range("R3").Value = 3477.043778601615
$ws.Range("S3").Value = 0.07676477306436616
$ws.Range("T3").Value = 0.07676477306436617
$ws.Range("G4").Value = 35.71508866666667
$ws.Range("H4").Value = 107.145266
$ws.Range("I4").Value = 0.1390302752364672
$ws.Range("J4").Value = 0.1390302752364672
$ws.Range("M4").Value = 1.973864333333333
$ws.Range("N4").Value = 5.921593
$ws.Range("O4").Value = 0.1007520683922509
$ws.Range("P4").Value = 0.1007520683922509
$ws.Range("Q4").Value = 70.4967396809709
$ws.Range("R4").Value = 634.4706571287381
$ws.Range("S4").Value = 0.01400758779921801
$ws.Range("T4").Value = 0.01400758779921801
$ws.Range("G5").Value = 35.71508866666667
$ws.Range("H5").Value = 107.145266
$ws.Range("I5").Value = 0.1390302752364672
$ws.Range("J5").Value = 0.1390302752364672
$ws.Range("M5").Value = 5.821642
$ws.Range("N5").Value = 17.464926
$ws.Range("O5").Value = 0.2971544006515816
$ws.Range("P5").Value = 0.2971544006515816
$ws.Range("Q5").Value = 207.9204602155907
$ws.Range("R5").Value = 1871.284141940316
$ws.Range("S5").Value = 0.04131345811031683
$ws.Range("T5").Value = 0.04131345811031684
$ws.Range("G6").Value = 54.09018966666667
$ws.Range("I6").Value = 0.2105601368412127
$ws.Range("J6").Value = 0.2105601368412127
$ws.Range("M6").Value = 0.9785706666666668
$ws.Range("N6").Value = 2.935712000000001
$ws.Range("O6").Value = 0.04994923768046061
$ws.Range("P6").Value = 0.04994923768046061
$ws.Range("Q6").Value = 52.93107296223646
$ws.Range("R6").Value = 476.3796566601281
$ws.Range("S6").Value = 0.01051731832111204
$ws.Range("T6").Value = 0.01051731832111204
$ws.Range("G7").Value = 54.09018966666667
$ws.Range("I7").Value = 0.2105601368412127
$ws.Range("J7").Value = 0.2105601368412127
$ws.Range("O7").Value = 0.5521442932757068
$ws.Range("P7").Value = 0.5521442932757068
$ws.Range("Q7").Value = 585.1058240372612
$ws.Range("R7").Value = 5265.952416335351
$ws.Range("S7").Value = 0.1162595779482275
$ws.Range("T7").Value = 0.1162595779482275
$ws.Range("G8").Value = 54.09018966666667
$ws.Range("I8").Value = 0.2105601368412127
$ws.Range("J8").Value = 0.2105601368412127
$ws.Range("M8").Value = 1.973864333333333
$ws.Range("N8").Value = 5.921593
$ws.Range("O8").Value = 0.1007520683922509
$ws.Range("P8").Value = 0.1007520683922509
$ws.Range("Q8").Value = 106.7666961662685
$ws.Range("R8").Value = 960.9002654964169
$ws.Range("S8").Value = 0.02121436930770757
$ws.Range("T8").Value = 0.02121436930770757
$ws.Range("G9").Value = 54.09018966666667
$ws.Range("I9").Value = 0.2105601368412127
$ws.Range("J9").Value = 0.2105601368412127
$ws.Range("M9").Value = 5.821642
$ws.Range("N9").Value = 17.464926
$ws.Range("O9").Value = 0.2971544006515816
$ws.Range("P9").Value = 0.2971544006515816
$ws.Range("Q9").Value = 314.8937199514327
$ws.Range("R9").Value = 2834.043479562893
$ws.Range("S9").Value = 0.06256887126416556
$ws.Range("T9").Value = 0.06256887126416556
$ws.Range("G10").Value = 101.4529346666666
$ws.Range("H10").Value = 304.358804
$ws.Range("I10").Value = 0.3949319449238378
$ws.Range("J10").Value = 0.3949319449238378
$ws.Range("M10").Value = 0.9785706666666668
$ws.Range("N10").Value = 2.935712000000001
$ws.Range("O10").Value = 0.04994923768046061
$ws.Range("P10").Value = 0.04994923768046061
$ws.Range("Q10").Value = 99.27886591204978
$ws.Range("R10").Value = 893.5097932084481
$ws.Range("S10").Value = 0.01972654958460735
$ws.Range("T10").Value = 0.01972654958460735
$ws.Range("G11").Value = 101.4529346666666
$ws.Range("H11").Value = 304.358804
$ws.Range("I11").Value = 0.3949319449238378
$ws.Range("J11").Value = 0.3949319449238378
$ws.Range("O11").Value = 0.5521442932757068
$ws.Range("P11").Value = 0.5521442932757068
$ws.Range("Q11").Value = 1097.439356470212
$ws.Range("R11").Value = 9876.954208231915
$ws.Range("S11").Value = 0.2180594196219728
$ws.Range("T11").Value = 0.2180594196219728
$ws.Range("G12").Value = 101.4529346666666
$ws.Range("H12").Value = 304.358804
$ws.Range("I12").Value = 0.3949319449238378
$ws.Range("J12").Value = 0.3949319449238378
$ws.Range("M12").Value = 1.973864333333333
$ws.Range("N12").Value = 5.921593
$ws.Range("O12").Value = 0.1007520683922509
$ws.Range("P12").Value = 0.1007520683922509
$ws.Range("Q12").Value = 200.2543292505302
$ws.Range("R12").Value = 1802.288963254772
$ws.Range("S12").Value = 0.03979021032525117
$ws.Range("T12").Value = 0.03979021032525118
$ws.Range("G13").Value = 101.4529346666666
$ws.Range("H13").Value = 304.358804
$ws.Range("I13").Value = 0.3949319449238378
$ws.Range("J13").Value = 0.3949319449238378
$ws.Range("M13").Value = 5.821642
$ws.Range("N13").Value = 17.464926
$ws.Range("O13").Value = 0.2971544006515816
$ws.Range("P13").Value = 0.2971544006515816
$ws.Range("Q13").Value = 590.6226654787225
$ws.Range("R13").Value = 5315.603989308503
$ws.Range("S13").Value = 0.1173557653920064
$ws.Range("T13").Value = 0.1173557653920065
$ws.Range("G14").Value = 65.628919
$ws.Range("H14").Value = 196.886757
$ws.Range("I14").Value = 0.2554776429984823
$ws.Range("J14").Value = 0.2554776429984823
$ws.Range("M14").Value = 0.9785706666666668
$ws.Range("N14").Value = 2.935712000000001
$ws.Range("O14").Value = 0.04994923768046061
$ws.Range("P14").Value = 0.04994923768046061
$ws.Range("Q14").Value = 64.22253501844267
$ws.Range("R14").Value = 578.0028151659841
$ws.Range("S14").Value = 0.01276091351217506
$ws.Range("T14").Value = 0.01276091351217505
$ws.Range("G15").Value = 65.628919
$ws.Range("H15").Value = 196.886757
$ws.Range("I15").Value = 0.2554776429984823
$ws.Range("J15").Value = 0.2554776429984823
$ws.Range("O15").Value = 0.5521442932757068
$ws.Range("P15").Value = 0.5521442932757068
$ws.Range("Q15").Value = 709.9228708350003
$ws.Range("R15").Value = 6389.305837515002
$ws.Range("S15").Value = 0.1410605226411403
$ws.Range("T15").Value = 0.1410605226411403
$ws.Range("G16").Value = 65.628919
$ws.Range("H16").Value = 196.886757
$ws.Range("I16").Value = 0.2554776429984823
$ws.Range("J16").Value = 0.2554776429984823
$ws.Range("M16").Value = 1.973864333333333
$ws.Range("N16").Value = 5.921593
$ws.Range("O16").Value = 0.1007520683922509
$ws.Range("P16").Value = 0.1007520683922509
$ws.Range("Q16").Value = 129.5425824493223
$ws.Range("R16").Value = 1165.883242043901
$ws.Range("S16").Value = 0.02573990096007415
$ws.Range("T16").Value = 0.02573990096007415
$ws.Range("G17").Value = 65.628919
$ws.Range("H17").Value = 196.886757
$ws.Range("I17").Value = 0.2554776429984823
$ws.Range("J17").Value = 0.2554776429984823
$ws.Range("M17").Value = 5.821642
$ws.Range("N17").Value = 17.464926
$ws.Range("O17").Value = 0.2971544006515816
$ws.Range("P17").Value = 0.2971544006515816
$ws.Range("Q17").Value = 382.068071264998
$ws.Range("R17").Value = 3438.612641384982
$ws.Range("S17").Value = 0.07591630588509274
$ws.Range("T17").Value = 0.07591630588509274
